$wb = $excel.ActiveWorkbook

$wsRun = $wb.Worksheets.Item("run_settings")
$wsFiles = $wb.Worksheets.Item("file_locations")
$wsVars = $wb.Worksheets.Item("variables")

# Update values on run_settings sheet
$wsRun.Range("B3").Value = 94

$falseCells = @("B19", "B20", "B24", "B27", "B29")
foreach ($addr in $falseCells) {
    $rng = $wsRun.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = "FALSE"
}

$wb.Save()
